# Updated cryptos list (price + 1h volume change) refresh.
# Price (col D) and Volume(1h) (col E) are stored as plain text in the sheet,
# e.g. "59.456.30" / "  -1.89%  " (double-dot thousands separators and
# padded percentages aren't valid Excel numbers). Where the new price text
# *would* parse as a real number (e.g. "560.93"), force the cell to Text
# format first so it is written back as a literal string, matching the
# original inline-string cell type instead of silently becoming numeric.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.459.79"
$ws.Range("E2").Value = "  -2.49%  "

$ws.Range("D3").Value = "2.589.18"
$ws.Range("E3").Value = "  -2.26%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.93"
$ws.Range("E5").Value = "  -1.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.37"
$ws.Range("E6").Value = "  -2.50%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("E8").Value = "  -2.05%  "

$ws.Range("D9").Value = "2.598.83"
$ws.Range("E9").Value = "  -2.89%  "

$ws.Range("E10").Value = "  -2.58%  "

$ws.Range("E11").Value = "  -0.69%  "

$ws.Range("E12").Value = "  +10.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.357"
$ws.Range("E13").Value = "  +4.04%  "

$ws.Range("D14").Value = "3.047.53"
$ws.Range("E14").Value = "  -2.16%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.38"
$ws.Range("E15").Value = "  +6.87%  "

$ws.Range("D16").Value = "59.393.70"
$ws.Range("E16").Value = "  -2.41%  "

$ws.Range("E17").Value = "  -0.10%  "

$ws.Range("D18").Value = "2.592.14"
$ws.Range("E18").Value = "  -2.62%  "

$ws.Range("E19").Value = "  +0.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "338.41"
$ws.Range("E20").Value = "  -1.78%  "

$ws.Range("E21").Value = "  -0.67%  "

$ws.Range("E22").Value = "  +1.35%  "

$ws.Range("E23").Value = "  +0.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.82"
$ws.Range("E24").Value = "  -4.24%  "

$ws.Range("E25").Value = "  +6.12%  "

$ws.Range("E26").Value = "  +0.49%  "

$ws.Range("E27").Value = "  -2.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.44"
$ws.Range("E28").Value = "  +0.21%  "

$ws.Range("E29").Value = "  -0.74%  "

$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("E31").Value = "  -3.70%  "

$ws.Range("E32").Value = "  -2.84%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.25"
$ws.Range("E33").Value = "  +2.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.08"
$ws.Range("E34").Value = "  -0.87%  "

$ws.Range("E35").Value = "  -0.90%  "

$ws.Range("E36").Value = "  -1.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.886"
$ws.Range("E37").Value = "  -2.44%  "

$ws.Range("E38").Value = "  -3.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.45"
$ws.Range("E39").Value = "  -0.36%  "

$ws.Range("E40").Value = "  -2.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.68"
$ws.Range("E41").Value = "  +0.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "294.12"
$ws.Range("E42").Value = "  -3.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "136.38"
$ws.Range("E43").Value = "  +6.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  +0.17%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0976"
$ws.Range("E45").Value = "  -0.90%  "

$ws.Range("E46").Value = "  -1.77%  "

$ws.Range("E47").Value = "  -2.64%  "

$ws.Range("E48").Value = "  -0.21%  "

$ws.Range("E49").Value = "  -0.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.81"
$ws.Range("E50").Value = "  -1.25%  "

$ws.Range("D51").Value = "1.956.46"
$ws.Range("E51").Value = "  -0.66%  "
